$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume (E) columns keep their original text formatting
# so numeric-looking strings are not coerced into floating point numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range('D2').Value = '70.860.81'
$ws.Range('E2').Value = '  +0.43%  '

# Row 3
$ws.Range('D3').Value = '3.537.78'
$ws.Range('E3').Value = '  -0.73%  '

# Row 4
$ws.Range('E4').Value = '  +0.03%  '

# Row 5
$ws.Range('D5').Value = '615.21'
$ws.Range('E5').Value = '  +0.02%  '

# Row 6
$ws.Range('D6').Value = '173.62'
$ws.Range('E6').Value = '  +0.74%  '

# Row 7
$ws.Range('B7').Value = 'LidoStakedEther'
$ws.Range('C7').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D7').Value = '3.532.63'
$ws.Range('E7').Value = '  -0.84%  '

# Row 8
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').Value = '0.612'
$ws.Range('E8').Value = '  -0.98%  '

# Row 9
$ws.Range('E9').Value = '  +0.06%  '

# Row 10
$ws.Range('E10').Value = '  +0.91%  '

# Row 11
$ws.Range('E11').Value = '  +0.46%  '

# Row 12
$ws.Range('D12').Value = '0.589'
$ws.Range('E12').Value = '  +0.61%  '

# Row 13
$ws.Range('D13').Value = '46.67'
$ws.Range('E13').Value = '  +0.25%  '

# Row 14
$ws.Range('E14').Value = '  +0.07%  '

# Row 15
$ws.Range('D15').Value = '4.107.87'
$ws.Range('E15').Value = '  -0.73%  '

# Row 16
$ws.Range('E16').Value = '  +1.14%  '

# Row 17
$ws.Range('D17').Value = '612.62'
$ws.Range('E17').Value = '  -0.48%  '

# Row 18
$ws.Range('D18').Value = '3.533.74'
$ws.Range('E18').Value = '  -1.23%  '

# Row 19
$ws.Range('D19').Value = '70.891.47'
$ws.Range('E19').Value = '  +0.37%  '

# Row 20
$ws.Range('E20').Value = '  +1.45%  '

# Row 21
$ws.Range('D21').Value = '17.81'
$ws.Range('E21').Value = '  +2.29%  '

# Row 22
$ws.Range('D22').Value = '0.889'
$ws.Range('E22').Value = '  +0.94%  '

# Row 23
$ws.Range('D23').Value = '9.04'
$ws.Range('E23').Value = '  -3.77%  '

# Row 24
$ws.Range('D24').Value = '15.75'
$ws.Range('E24').Value = '  +0.03%  '

# Row 25
$ws.Range('D25').Value = '98.32'
$ws.Range('E25').Value = '  +1.71%  '

# Row 26
$ws.Range('E26').Value = '  -1.28%  '

# Row 27
$ws.Range('E27').Value = '  +0.08%  '

# Row 28
$ws.Range('E28').Value = '  -0.23%  '

# Row 29
$ws.Range('D29').Value = '33.92'
$ws.Range('E29').Value = '  +1.14%  '

# Row 30
$ws.Range('D30').Value = '9.17'
$ws.Range('E30').Value = '  +1.47%  '

# Row 31
$ws.Range('D31').Value = '3.03'
$ws.Range('E31').Value = '  -0.40%  '

# Row 32
$ws.Range('D32').Value = '8.19'
$ws.Range('E32').Value = '  -3.57%  '

# Row 33
$ws.Range('D33').Value = '1.30'
$ws.Range('E33').Value = '  -0.03%  '

# Row 34
$ws.Range('D34').Value = '6.88'
$ws.Range('E34').Value = '  -0.91%  '

# Row 35
$ws.Range('D35').Value = '605.44'
$ws.Range('E35').Value = '  +5.59%  '

# Row 36
$ws.Range('E36').Value = '  -0.64%  '

# Row 37
$ws.Range('E37').Value = '  +0.45%  '

# Row 38
$ws.Range('E38').Value = '  -2.64%  '

# Row 39
$ws.Range('D39').Value = '0.0472'
$ws.Range('E39').Value = '  +0.76%  '

# Row 40
$ws.Range('D40').Value = '57.01'
$ws.Range('E40').Value = '  -0.53%  '

# Row 41
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '0.146'
$ws.Range('E41').Value = '  +2.99%  '

# Row 42
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.14%  '

# Row 43
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '3.387.58'
$ws.Range('E43').Value = '  +0.17%  '

# Row 44
$ws.Range('B44').Value = 'PEPE'
$ws.Range('C44').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D44').Value = '0.0₃0742'
$ws.Range('E44').Value = '  +5.89%  '

# Row 45
$ws.Range('D45').Value = '2.98'
$ws.Range('E45').Value = '  +0.51%  '

# Row 46
$ws.Range('E46').Value = '  -1.69%  '

# Row 47
$ws.Range('D47').Value = '32.28'
$ws.Range('E47').Value = '  -2.10%  '

# Row 48
$ws.Range('E48').Value = '  -1.47%  '

# Row 49
$ws.Range('E49').Value = '  +0.78%  '

# Row 50
$ws.Range('D50').Value = '133.79'
$ws.Range('E50').Value = '  +0.11%  '

# Row 51
$ws.Range('E51').Value = '  -0.01%  '
